$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ("Ready for handoff" -eq $val) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Narrow the "Status" columns to fit the new (shorter) text ---
# Overview sheet: Status is mirrored in columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns("E:F").ColumnWidth = 12.5

# zh-cn / de-de sheets: Status is column C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns("C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns("C").ColumnWidth = 12.5
